# Auto-generated edit script: updates FFXIV market-price derived cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 41.625
$ws.Range("I39").Value = 41.625
$ws.Range("K39").Value = 124.875
$ws.Range("M39").Value = 171.125

$ws.Range("H41").Value = 622.5357
$ws.Range("I41").Value = 503.6
$ws.Range("J41").Value = 759.7692
$ws.Range("K41").Value = 503.6
$ws.Range("L41").Value = 759.7692
$ws.Range("M41").Value = -63.60000000000002
$ws.Range("N41").Value = -1639.7692

$ws.Range("H106").Value = 8525.058999999999
$ws.Range("I106").Value = 2106.889
$ws.Range("K106").Value = 2106.889
$ws.Range("M106").Value = -1475.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 23926.285
$ws.Range("J24").Value = 23926.285
$ws.Range("L24").Value = 23926.285
$ws.Range("N24").Value = -24674.285

$ws.Range("H63").Value = 2363
$ws.Range("I63").Value = 2359.3
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 2359.3
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -1673.3
$ws.Range("N63").Value = -3772

$ws.Range("H66").Value = 2363
$ws.Range("I66").Value = 2359.3
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 11796.5
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -8364.5
$ws.Range("N66").Value = -18864

$ws.Range("H74").Value = 2372.0588
$ws.Range("I74").Value = 2021.9333
$ws.Range("J74").Value = 4998
$ws.Range("K74").Value = 2021.9333
$ws.Range("L74").Value = 4998
$ws.Range("M74").Value = -1147.9333
$ws.Range("N74").Value = -6746

$ws.Range("H77").Value = 2372.0588
$ws.Range("I77").Value = 2021.9333
$ws.Range("J77").Value = 4998
$ws.Range("K77").Value = 10109.6665
$ws.Range("L77").Value = 24990
$ws.Range("M77").Value = -5741.666499999999
$ws.Range("N77").Value = -33726

$ws.Range("H100").Value = 23926.285
$ws.Range("J100").Value = 23926.285
$ws.Range("L100").Value = 23926.285
$ws.Range("N100").Value = -26090.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 30000
$ws.Range("J35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30620

$ws.Range("H94").Value = 754.25
$ws.Range("I94").Value = 735.6667
$ws.Range("J94").Value = 810
$ws.Range("K94").Value = 735.6667
$ws.Range("L94").Value = 810
$ws.Range("M94").Value = -284.6667
$ws.Range("N94").Value = -1712

$ws.Range("H134").Value = 969.4
$ws.Range("I134").Value = 988.7826
$ws.Range("K134").Value = 2966.3478
$ws.Range("M134").Value = -431.3478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1970.925
$ws.Range("I31").Value = 1574
$ws.Range("K31").Value = 1574
$ws.Range("M31").Value = -1279

$ws.Range("H34").Value = 1970.925
$ws.Range("I34").Value = 1574
$ws.Range("K34").Value = 1574
$ws.Range("M34").Value = -1372

$ws.Range("H134").Value = 4133
$ws.Range("I134").Value = 4133
$ws.Range("K134").Value = 12399
$ws.Range("M134").Value = -9864

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1122.25
$ws.Range("I5").Value = 266.33334
$ws.Range("K5").Value = 799.0000200000001
$ws.Range("M5").Value = -687.0000200000001

$ws.Range("H70").Value = 4163.4287
$ws.Range("I70").Value = 2786
$ws.Range("K70").Value = 8358
$ws.Range("M70").Value = -8043

$ws.Range("H73").Value = 4163.4287
$ws.Range("I73").Value = 2786
$ws.Range("K73").Value = 8358
$ws.Range("M73").Value = -7266

$ws.Range("H132").Value = 2648.3333
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 2878
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 25902
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -30962

$ws.Range("H134").Value = 1353.875
$ws.Range("I134").Value = 1353.875
$ws.Range("K134").Value = 4061.625
$ws.Range("M134").Value = 1008.375

$ws.Range("H135").Value = 1122.25
$ws.Range("I135").Value = 266.33334
$ws.Range("K135").Value = 2397.00006
$ws.Range("M135").Value = 137.9999399999997

$ws.Range("H139").Value = 2539
$ws.Range("I139").Value = 2779.923
$ws.Range("J139").Value = 1495
$ws.Range("K139").Value = 8339.769
$ws.Range("L139").Value = 4485
$ws.Range("M139").Value = -3199.769
$ws.Range("N139").Value = -14765

$ws.Range("H140").Value = 4217.25
$ws.Range("I140").Value = 1680.0625
$ws.Range("J140").Value = 9291.625
$ws.Range("K140").Value = 5040.1875
$ws.Range("L140").Value = 27874.875
$ws.Range("M140").Value = 139.8125
$ws.Range("N140").Value = -38234.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 37560.75
$ws.Range("J95").Value = 37560.75
$ws.Range("L95").Value = 37560.75
$ws.Range("N95").Value = -43052.75

$ws.Range("H132").Value = 5667.7144
$ws.Range("I132").Value = 5667.7144
$ws.Range("K132").Value = 17003.1432
$ws.Range("M132").Value = -14473.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 824.5
$ws.Range("I22").Value = 483.33334
$ws.Range("J22").Value = 1165.6666
$ws.Range("K22").Value = 483.33334
$ws.Range("L22").Value = 1165.6666
$ws.Range("M22").Value = -188.33334
$ws.Range("N22").Value = -1755.6666

$ws.Range("H27").Value = 824.5
$ws.Range("I27").Value = 483.33334
$ws.Range("J27").Value = 1165.6666
$ws.Range("K27").Value = 483.33334
$ws.Range("L27").Value = 1165.6666
$ws.Range("M27").Value = -376.33334
$ws.Range("N27").Value = -1379.6666

$ws.Range("H107").Value = 3999.5
$ws.Range("I107").Value = 3999.5
$ws.Range("K107").Value = 3999.5
$ws.Range("M107").Value = -2079.5

$ws.Range("H132").Value = 8676.888999999999
$ws.Range("I132").Value = 11028.182
$ws.Range("K132").Value = 33084.546
$ws.Range("M132").Value = -30554.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws.Range("H33").Value = 27500
$ws.Range("J33").Value = 27500
$ws.Range("L33").Value = 27500
$ws.Range("N33").Value = -28000

$ws.Range("H36").Value = 27500
$ws.Range("J36").Value = 27500
$ws.Range("L36").Value = 27500
$ws.Range("N36").Value = -28000

$ws.Range("H81").Value = 2577.25
$ws.Range("I81").Value = 1915.8334
$ws.Range("K81").Value = 3831.6668
$ws.Range("M81").Value = -2770.6668

$ws.Range("H84").Value = 2577.25
$ws.Range("I84").Value = 1915.8334
$ws.Range("K84").Value = 19158.334
$ws.Range("M84").Value = -13854.334

$ws.Range("H132").Value = 13325.218
$ws.Range("J132").Value = 14345
$ws.Range("L132").Value = 43035
$ws.Range("N132").Value = -48095
